$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "https://www.syngenta.bg/products/search/crop-protection"
$ws.Range("A6").Select()
